$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.949.04"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.92%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.817.63"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.29%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.04%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.65"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.54%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.61"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.42%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.816.64"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -0.92%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +0.58%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +9.38%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.98"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.58%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.459.62"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -2.38%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.832.61"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -2.19%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.64"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +2.45%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.045.22"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -0.79%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  +0.07%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.88"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.19%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.09"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.45%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.741"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -9.03%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.66"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  +2.40%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.03%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("D28").Style = $style
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -1.56%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.966.26"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("E33").Value = "  -1.22%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.83"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -2.35%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.38"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -0.80%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.784.62"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  +4.26%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.99"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +0.74%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.140"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("E41").Value = "  -2.38%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("E44").Value = "  +2.15%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -1.60%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "411.73"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.51"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000286"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -4.30%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.63"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -0.73%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0360"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.35%  "
